# Append new scrape results (2026-01-22 12:43 JST) to the top of the
# "ランサーズ" sheet, pushing the existing rows down, and widen column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove existing hyperlinks first; their ranges will be stale once we
#     shift rows down, and we'll re-add fresh ones at the end in the right
#     order (F2..F8).
$ws.UsedRange.Hyperlinks.Delete()

# --- Insert 4 new blank rows above the existing row 2, shifting the
#     3 existing data rows from 2-4 down to 6-8.
$ws.Range("A2:A5").EntireRow.Insert(-4121)

# --- Update the timestamp on the rows that already existed (now at 6-8)
$ws.Range("A6").Value = "2026-01-22 12:43:52"
$ws.Range("A7").Value = "2026-01-22 12:43:52"
$ws.Range("A8").Value = "2026-01-22 12:43:52"

# --- New row 2: 産業機械向けAI異常検知...
$ws.Range("A2").Value = "2026-01-22 12:43:52"
$ws.Range("B2").Value = "産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5450864"
$ws.Range("G2").Value = 383
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- New row 3: 【フルタイム】最先端AI(LLM)...
$ws.Range("A3").Value = "2026-01-22 12:43:52"
$ws.Range("B3").Value = "【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5460294"
$ws.Range("G3").Value = 375
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# --- New row 4: 【週5日】法人向け生成AIサービス...
$ws.Range("A4").Value = "2026-01-22 12:43:52"
$ws.Range("B4").Value = "【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5460267"
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"

# --- New row 5: 初回 Keepaの通知から...
$ws.Range("A5").Value = "2026-01-22 12:43:52"
$ws.Range("B5").Value = "初回 Keepaの通知からAmazonで自動購入するシステムの開発依頼の仕事"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5476963"
$ws.Range("G5").Value = 75
$ws.Range("H5").Value = "◆開発"

# --- Re-create the hyperlinks for F2:F8 in order, matching rId1..rId7
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5450864")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5460294")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5460267")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5476963")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5476347")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5476708")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5476581")

# --- Widen column B from 38 to 51 characters
$ws.Columns.Item(2).ColumnWidth = 50.15
